$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The decision table previously had a duplicated third condition column (D)
# for rows 13-16, mirroring column C. Remove the now-unneeded column D
# values so that only two condition columns (B, C) remain for this table,
# matching the new behavior where columns with true conditions are matched
# with the output object once the number of IsTrue conditions exceeds 8.
$ws.Range("D13:D16").ClearContents()

# Refresh the active selection to reflect where the author was working
# when the change was made.
[void]$ws.Range("H19").Select()
